$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Swap FP (E) and FN (F) counts for the rows whose values changed:
# Mediastore (79), TEAMMATES_H (81), Teastore_H (85), JabRef_H (87)
foreach ($r in 79, 81, 85, 87) {
    $e = $ws.Cells.Item($r, 5).Value()
    $f = $ws.Cells.Item($r, 6).Value()
    $ws.Cells.Item($r, 5).Value = $f
    $ws.Cells.Item($r, 6).Value = $e
}

# Row 82 (TEAMMATES) was missing the Accuracy/Specificity formulas that
# every other project row has - add them now. Re-entering the same
# formula across the whole K80:K87 / L80:L87 block (rather than touching
# only K82/L82) lets Excel keep using its shared-formula grouping for the
# block instead of emitting a new one-off formula just for row 82.
$ws.Range("K80:K87").Formula = "=(D80+G80)/SUM(D80:G80)"
$ws.Range("L80:L87").Formula = "=G80/(G80+E80)"

# Update the view state (scroll position + active selection) to match.
$ws.Application.ActiveWindow.ScrollRow = 71
$ws.Range("R82").Select()
